$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1175.8334
$ws.Range("J103").Value = 1278.0667
$ws.Range("L103").Value = 3834.2001
$ws.Range("N103").Value = -5006.2001
$ws.Range("H125").Value = 100001960
$ws.Range("J125").Value = 2450.5
$ws.Range("L125").Value = 22054.5
$ws.Range("N125").Value = -26974.5
$ws.Range("H134").Value = 80779.5
$ws.Range("J134").Value = 80779.5
$ws.Range("L134").Value = 80779.5
$ws.Range("N134").Value = -90919.5
$ws.Range("H137").Value = 3028.9473
$ws.Range("J137").Value = 2573.875
$ws.Range("L137").Value = 7721.625
$ws.Range("N137").Value = -12821.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3299.8823
$ws.Range("I2").Value = 908.2727
$ws.Range("J2").Value = 7684.5
$ws.Range("K2").Value = 908.2727
$ws.Range("L2").Value = 7684.5
$ws.Range("M2").Value = -795.2727
$ws.Range("N2").Value = -7910.5
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184
$ws.Range("H32").Value = 1510034.4
$ws.Range("I32").Value = 1648176.4
$ws.Range("K32").Value = 1648176.4
$ws.Range("M32").Value = -1647889.4
$ws.Range("H61").Value = 6623.875
$ws.Range("I61").Value = 3238.2222
$ws.Range("J61").Value = 13655.615
$ws.Range("K61").Value = 3238.2222
$ws.Range("L61").Value = 13655.615
$ws.Range("M61").Value = -3026.2222
$ws.Range("N61").Value = -14079.615
$ws.Range("H74").Value = 34124.527
$ws.Range("I74").Value = 44169
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 44169
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -43295
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 34124.527
$ws.Range("I77").Value = 44169
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 220845
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -216477
$ws.Range("N77").Value = -38736
$ws.Range("H110").Value = 14493637
$ws.Range("I110").Value = 915.55
$ws.Range("J110").Value = 111111784
$ws.Range("K110").Value = 915.55
$ws.Range("L110").Value = 111111784
$ws.Range("M110").Value = 1129.45
$ws.Range("N110").Value = -111115874
$ws.Range("H116").Value = 3299.8823
$ws.Range("I116").Value = 908.2727
$ws.Range("J116").Value = 7684.5
$ws.Range("K116").Value = 908.2727
$ws.Range("L116").Value = 7684.5
$ws.Range("M116").Value = 1385.7273
$ws.Range("N116").Value = -12272.5
$ws.Range("H122").Value = 46800
$ws.Range("J122").Value = 10666.667
$ws.Range("L122").Value = 32000.001
$ws.Range("N122").Value = -36900.001
$ws.Range("H136").Value = 6623.875
$ws.Range("I136").Value = 3238.2222
$ws.Range("J136").Value = 13655.615
$ws.Range("K136").Value = 9714.6666
$ws.Range("L136").Value = 40966.845
$ws.Range("M136").Value = -7164.6666
$ws.Range("N136").Value = -46066.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3299.8823
$ws.Range("I3").Value = 908.2727
$ws.Range("J3").Value = 7684.5
$ws.Range("K3").Value = 908.2727
$ws.Range("L3").Value = 7684.5
$ws.Range("M3").Value = -794.2727
$ws.Range("N3").Value = -7912.5
$ws.Range("H86").Value = 37076516
$ws.Range("I86").Value = 54489
$ws.Range("K86").Value = 54489
$ws.Range("M86").Value = -53366
$ws.Range("H89").Value = 37076516
$ws.Range("I89").Value = 54489
$ws.Range("K89").Value = 272445
$ws.Range("M89").Value = -266829
$ws.Range("H94").Value = 3329.0908
$ws.Range("I94").Value = 929
$ws.Range("J94").Value = 5329.1665
$ws.Range("K94").Value = 929
$ws.Range("L94").Value = 5329.1665
$ws.Range("M94").Value = -478
$ws.Range("N94").Value = -6231.1665
$ws.Range("H99").Value = 3499184
$ws.Range("I99").Value = 2672.5
$ws.Range("J99").Value = 11366335
$ws.Range("K99").Value = 2672.5
$ws.Range("L99").Value = 11366335
$ws.Range("M99").Value = -1174.5
$ws.Range("N99").Value = -11369331
$ws.Range("H107").Value = 41670064
$ws.Range("I107").Value = 56252724
$ws.Range("J107").Value = 5318.143
$ws.Range("K107").Value = 56252724
$ws.Range("L107").Value = 5318.143
$ws.Range("M107").Value = -56250804
$ws.Range("N107").Value = -9158.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6850.2354
$ws.Range("I31").Value = 2856.76
$ws.Range("K31").Value = 2856.76
$ws.Range("M31").Value = -2561.76
$ws.Range("H34").Value = 6850.2354
$ws.Range("I34").Value = 2856.76
$ws.Range("K34").Value = 2856.76
$ws.Range("M34").Value = -2654.76
$ws.Range("H58").Value = 10874929
$ws.Range("I58").Value = 16668210
$ws.Range("J58").Value = 12527.4375
$ws.Range("K58").Value = 16668210
$ws.Range("L58").Value = 12527.4375
$ws.Range("M58").Value = -16668007
$ws.Range("N58").Value = -12933.4375
$ws.Range("H60").Value = 19856.428
$ws.Range("I60").Value = 12800
$ws.Range("J60").Value = 37497.5
$ws.Range("K60").Value = 12800
$ws.Range("L60").Value = 37497.5
$ws.Range("M60").Value = -12289
$ws.Range("N60").Value = -38519.5
$ws.Range("H96").Value = 26983
$ws.Range("J96").Value = 26983
$ws.Range("L96").Value = 26983
$ws.Range("N96").Value = -32475
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H136").Value = 10874929
$ws.Range("I136").Value = 16668210
$ws.Range("J136").Value = 12527.4375
$ws.Range("K136").Value = 50004630
$ws.Range("L136").Value = 37582.3125
$ws.Range("M136").Value = -50002080
$ws.Range("N136").Value = -42682.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1572686.5
$ws.Range("J122").Value = 1236.5
$ws.Range("L122").Value = 11128.5
$ws.Range("N122").Value = -16028.5
$ws.Range("H132").Value = 7631.591
$ws.Range("I132").Value = 3642.8572
$ws.Range("J132").Value = 9493
$ws.Range("K132").Value = 32785.7148
$ws.Range("L132").Value = 85437
$ws.Range("M132").Value = -30255.7148
$ws.Range("N132").Value = -90497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 924.3333
$ws.Range("I97").Value = 856.17645
$ws.Range("K97").Value = 856.17645
$ws.Range("M97").Value = -360.17645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 41749.668
$ws.Range("J62").Value = 41749.668
$ws.Range("L62").Value = 41749.668
$ws.Range("N62").Value = -42997.668
$ws.Range("H65").Value = 41749.668
$ws.Range("J65").Value = 41749.668
$ws.Range("L65").Value = 125249.004
$ws.Range("N65").Value = -131489.004
$ws.Range("H93").Value = 7532.4546
$ws.Range("I93").Value = 3991.4
$ws.Range("K93").Value = 3991.4
$ws.Range("M93").Value = -2743.4
$ws.Range("H136").Value = 11308.218
$ws.Range("I136").Value = 2848.3333
$ws.Range("K136").Value = 8544.999899999999
$ws.Range("M136").Value = -5994.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 13889631
$ws.Range("I107").Value = 544.875
$ws.Range("J107").Value = 41667804
$ws.Range("K107").Value = 1634.625
$ws.Range("L107").Value = 125003412
$ws.Range("M107").Value = 285.375
$ws.Range("N107").Value = -125007252
$ws.Range("H125").Value = 125000
$ws.Range("J125").Value = 125000
$ws.Range("L125").Value = 125000
$ws.Range("N125").Value = -134840
$ws.Range("H132").Value = 15164632
$ws.Range("I132").Value = 25006266
$ws.Range("K132").Value = 75018798
$ws.Range("M132").Value = -75016268
$ws.Range("H136").Value = 25671758
$ws.Range("I136").Value = 52632624
$ws.Range("K136").Value = 157897872
$ws.Range("M136").Value = -157895322
